# "Added some more data"
# - add a new "added" date column (E) to the existing dataset catalogue
# - append 10 new catalogue rows (76-85) describing newly added Kaggle datasets
# - widen column B a bit to fit the longer description text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New column E: header + "added" dates for all existing rows (2-75)
# ---------------------------------------------------------------------------

$ws.Range("E1").Value = "added"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)   # xlPasteFormats - copy D1's style onto E1

for ($r = 2; $r -le 75; $r++) {
    $ws.Cells.Item($r, 5).Value = 20190201
}

# ---------------------------------------------------------------------------
# 2. Ten new rows describing newly added datasets
# ---------------------------------------------------------------------------

$newRows = @(
    @("Top 5000 Youtube channels data from Socialblade.", "General metrics of top 5000 YouTube channels by Socialblade", "top-5000-youtube", "https://www.kaggle.com/mdhrumil/top-5000-youtube-channels-data-from-socialblade"),
    @("Countries of the World", "Country names linked to region, population, area size, GDP, mortality and more", "countries of the world", "https://www.kaggle.com/fernandol/countries-of-the-world"),
    @("Toy Dataset", "A dataset to play around with!", "toy_dataset", "https://www.kaggle.com/carlolepelaars/toy-dataset"),
    @("New York City Fee Charges", "From New York City Open Data", "fee-charges", "https://www.kaggle.com/new-york-city/new-york-city-fee-charges"),
    @("Google Play Store Apps", "Web scraped user review data of 10k Play Store apps for analysing the Android market.", "googleplaystore_user_reviews", "https://www.kaggle.com/lava18/google-play-store-apps"),
    @("Google Play Store Apps", "Web scraped data of 10k Play Store apps for analysing the Android market.", "googleplaystore", "https://www.kaggle.com/lava18/google-play-store-apps"),
    @("Video Game Sales with Ratings", "Video game sales from Vgchartz and corresponding ratings from Metacritic", "Video_Games_Sales_as_at_22_Dec_2016", "https://www.kaggle.com/rush4ratio/video-game-sales-with-ratings"),
    @("NBA Players stats since 1950", "Season data for NBA players since 1950", "Seasons_Stats", "https://www.kaggle.com/drgilermo/nba-players-stats"),
    @("NBA Players stats since 1950", "Generic player data for NBA players since 1950", "Players", "https://www.kaggle.com/drgilermo/nba-players-stats"),
    @("NBA Players stats since 1950", "Game related data for NBA players since 1950", "player_data", "https://www.kaggle.com/drgilermo/nba-players-stats")
)

$startRow = 76
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($row, 1).Value = $data[0]
    $ws.Cells.Item($row, 2).Value = $data[1]
    $ws.Cells.Item($row, 3).Value = $data[2]
    $ws.Cells.Item($row, 4).Value = $data[3]
    $ws.Cells.Item($row, 5).Value = 20190204
}

# ---------------------------------------------------------------------------
# 3. Widen column B so the longer descriptions fit
# ---------------------------------------------------------------------------

$ws.Columns.Item(2).ColumnWidth = 70.3

# ---------------------------------------------------------------------------
# 4. Update the on-screen selection to roughly match where editing left off
# ---------------------------------------------------------------------------

$ws.Range("B" + ($startRow + $newRows.Count + 1)).Select()
